# New weekly price record for Feria Lagunitas de Puerto Montt - Apio.
# Insert a new row at row 81 (pushes the former rows 81-121 down to 82-122,
# growing the used range from A1:R121 to A1:R122) and populate it with the
# new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(81).Insert()

$ws.Cells.Item(81, 1).Value2  = 4
$ws.Cells.Item(81, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(81, 3).Value2  = "Los Lagos"
$ws.Cells.Item(81, 4).Value2  = 44438
$ws.Cells.Item(81, 5).Value2  = 10
$ws.Cells.Item(81, 6).Value2  = 100112017
$ws.Cells.Item(81, 7).Value2  = "Apio"
$ws.Cells.Item(81, 8).Value2  = "Americana (o)"
$ws.Cells.Item(81, 9).Value2  = "Primera"
$ws.Cells.Item(81, 10).Value2 = 20
$ws.Cells.Item(81, 11).Value2 = 12000
$ws.Cells.Item(81, 12).Value2 = 12000
$ws.Cells.Item(81, 13).Value2 = 12000
$ws.Cells.Item(81, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(81, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(81, 16).Value2 = 2000
$ws.Cells.Item(81, 17).Value2 = 6
$ws.Cells.Item(81, 18).Value2 = "Hortaliza"
